$p = $ppt.ActivePresentation
$p.Slides.Item(1).Delete()
